$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (45181 -> 45182) for every data row (rows 2 through 91).
$ws.Range("C2:C91").Value = 45182
